$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 5 blank rows above the old summary row (old row 83) ---
# This pushes the old "totals" row (83) down to 88, and the trailing
# row (86) down to 91, picking up the same formatting as the
# surrounding blank rows 81/82.
$ws.Rows("83:87").Insert()

# --- 2. Fill in the new pub-run entry on row 81 (previously blank) ---
$ws.Range("A81").Value = 45672
$ws.Range("B81").Value = "Little Chester Ale House"
$ws.Range("C81").Value = "Little Chester"
$ws.Range("D81").Value = "start/end at pub"
$ws.Range("E81").Value = 4.96
$ws.Range("F81").Formula = "=E81*0.6213712"
$ws.Range("G81").Value = 0.022187499999999999
$ws.Range("H81").Formula = "=G81/F81"
$ws.Range("I81").Value = 1
$ws.Range("K81").Value = 1
$ws.Range("M81").Value = 1
$ws.Range("O81").Value = 1
$ws.Range("P81").Value = "Spooky mist"
$ws.Range("Q81").Formula = "=SUM(I81:O81)*F81"

# --- 3. Give the (previously empty) rows their running "Total Miles"
#        formula, matching the Q column's shared formula down the page ---
$ws.Range("Q80").Formula = "=SUM(I80:O80)*F80"
$ws.Range("Q82").Formula = "=SUM(I82:O82)*F82"
$ws.Range("Q83").Formula = "=SUM(I83:O83)*F83"
$ws.Range("Q84").Formula = "=SUM(I84:O84)*F84"
$ws.Range("Q85").Formula = "=SUM(I85:O85)*F85"
$ws.Range("Q86").Formula = "=SUM(I86:O86)*F86"

# The freshly inserted rows (83:87) don't inherit the sheet's small
# 8pt body font the way the pre-existing blank rows (81/82) do, so the
# new Q formulas above land on the default style instead of the
# shared "small font" cell style used everywhere else in the table.
# Nudge the font size to land back on that existing style.
$ws.Range("Q83:Q86").Font.Size = 8

# --- 4. Fix up the totals row, now shifted from row 83 to row 88, so
#        its ranges pick up the newly inserted rows (8:87 instead of
#        8:80 etc.) ---
$ws.Range("F88").Formula = "=SUM(F8:F87)"
$ws.Range("H88").Formula = "=AVERAGE(H9:H73)"
$ws.Range("I88").Formula = "=SUM(I3:I87)"
$ws.Range("J88").Formula = "=SUM(J3:J87)"
$ws.Range("K88").Formula = "=SUM(K3:K87)"
$ws.Range("L88").Formula = "=SUM(L3:L87)"
$ws.Range("M88").Formula = "=SUM(M3:M87)"
$ws.Range("N88").Formula = "=SUM(N3:N87)"
$ws.Range("O88").Formula = "=SUM(O3:O87)"
$ws.Range("Q88").Formula = "=SUM(Q3:Q87)"

# --- 5. Selection bookkeeping to match the saved workbook state ---
$ws.Range("Q89").Select()

